$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.533.63"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "1.810.63"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'225.66"
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("D6").Value = "'0.599"
$ws.Range("E6").Value = "  +2.64%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").Value = "'36.45"
$ws.Range("E8").Value = "  +4.80%  "

$ws.Range("E9").Value = "  -2.86%  "

$ws.Range("D10").Value = "'0.0683"
$ws.Range("E10").Value = "  -1.66%  "

$ws.Range("E11").Value = "  +1.58%  "

$ws.Range("D12").Value = "2.071.74"
$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("D13").Value = "'11.29"
$ws.Range("E13").Value = "  +0.72%  "

$ws.Range("D14").Value = "1.832.31"
$ws.Range("E14").Value = "  +1.67%  "

$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").Value = "34.500.54"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("D17").Value = "'4.43"
$ws.Range("E17").Value = "  +1.18%  "

$ws.Range("D18").Value = "'68.35"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").Value = "'242.70"
$ws.Range("E19").Value = "  -1.20%  "

$ws.Range("D20").Value = "0.0₃0776"
$ws.Range("E20").Value = "  -2.67%  "

$ws.Range("D21").Value = "'11.23"
$ws.Range("E21").Value = "  -2.38%  "

$ws.Range("E23").Value = "  -1.49%  "

$ws.Range("E24").Value = "  +5.47%  "

$ws.Range("D25").Value = "'171.64"
$ws.Range("E25").Value = "  -0.85%  "

$ws.Range("D26").Value = "'7.86"
$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").Value = "'17.27"
$ws.Range("E27").Value = "  +2.88%  "

$ws.Range("D28").Value = "'0.121"
$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("E30").Value = "  -0.30%  "

$ws.Range("E31").Value = "  -1.10%  "

$ws.Range("D32").Value = "'3.91"
$ws.Range("E32").Value = "  -2.28%  "

$ws.Range("D33").Value = "'0.0516"
$ws.Range("E33").Value = "  -2.57%  "

$ws.Range("D34").Value = "'1.80"
$ws.Range("E34").Value = "  -2.26%  "

$ws.Range("D35").Value = "1.362.86"
$ws.Range("E35").Value = "  -2.27%  "

$ws.Range("D36").Value = "'0.654"
$ws.Range("E36").Value = "  -3.80%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").Value = "'2.35"
$ws.Range("E38").Value = "  -6.06%  "

$ws.Range("E39").Value = "  -1.93%  "

$ws.Range("D40").Value = "'2.42"
$ws.Range("E40").Value = "  +0.72%  "

# Row 41/42: MXToken and Aave swapped positions, with new price/volume values
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'81.09"
$ws.Range("E41").Value = "  -2.77%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.78"
$ws.Range("E42").Value = "  -1.83%  "

$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("E44").Value = "  +3.99%  "

$ws.Range("D45").Value = "'13.43"
$ws.Range("E45").Value = "  -0.80%  "

$ws.Range("E46").Value = "  -2.56%  "

$ws.Range("D47").Value = "1.971.98"
$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("E48").Value = "  -2.67%  "

$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.29%  "

$ws.Range("D50").Value = "'102.59"
$ws.Range("E50").Value = "  -2.26%  "

$ws.Range("E51").Value = "  -6.11%  "
